# Apply "version final sin errores" edit:
#  - On the "Metadata" sheet, update Version value to 0.7.0
#  - Remove the "Jurisdiction" / "Chile" row entirely (row shifts up)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update the Version property value (row 3, column B)
$ws.Cells.Item(3, 2).Value = "0.7.0"

# Delete the entire "Jurisdiction" | "Chile" row (row 11), shifting subsequent rows up
$ws.Rows.Item(11).Delete()
